$wb = $excel.ActiveWorkbook

# user_assignments is the 4th sheet (system_permissions, system_roles, system_groups, user_assignments)
$ws = $wb.Worksheets.Item(4)

# --- Row 2: replace the old Super Admin email with the new one, keeping the
# hyperlink relationship id by deleting + re-adding the link on A2 ---------
$ws.Range("A2").Value = "rburgess@mpimedia.com"
$ws.Range("A2").Hyperlinks.Delete() | Out-Null
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:rburgess@mpimedia.com") | Out-Null
$ws.Range("A2").Style = "Hyperlink"
$ws.Range("B2").Value = "Super Admins"

# --- New row 3: dennis@dennmart.com ---------------------------------------
$ws.Range("A3").Value = "dennis@dennmart.com"
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:dennis@dennmart.com") | Out-Null
$ws.Range("A3").Style = "Hyperlink"
$ws.Range("B3").Value = "Super Admins"

# --- New row 4: badie@mpimedia.com ----------------------------------------
$ws.Range("A4").Value = "badie@mpimedia.com"
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:badie@mpimedia.com") | Out-Null
$ws.Range("A4").Style = "Hyperlink"
$ws.Range("B4").Value = "Super Admins"

# --- Make user_assignments the active/selected sheet, with B16 selected ---
$ws.Activate()
$ws.Range("B16").Select() | Out-Null

Write-Output "done"
